$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "K health_old" column (M) metrics to reflect the
# newly added Harvard case classification.
$ws.Range("M2").Value = 0.6666666666666666   # precision
$ws.Range("M3").Value = 1                     # recall
$ws.Range("M4").Value = 0.8                   # f1-score
$ws.Range("M5").Value = 0.9090909090909091    # f2-score
$ws.Range("M6").Value = 0.9639404333166532    # NDCG
